$d = $word.ActiveDocument

$replacements = @(
    @("2023-12-30 Saturday", "2023-12-31 Sunday"),
    @("70×28=1960", "50×63=3150"),
    @("66×14=924", "92×92=8464"),
    @("59×60=3540", "75×44=3300"),
    @("16×14=224", "54×73=3942"),
    @("35×88=3080", "99×33=3267"),
    @("65×54=3510", "29×81=2349"),
    @("61×23=1403", "21×34=714"),
    @("22×50=1100", "84×79=6636"),
    @("84×61=5124", "39×65=2535"),
    @("44×60=2640", "67×69=4623"),
    @("60×73=4380", "62×27=1674"),
    @("26×34=884", "99×57=5643"),
    @("96×12=1152", "24×20=480"),
    @("54×55=2970", "57×28=1596"),
    @("54×45=2430", "63×26=1638"),
    @("82×81=6642", "39×28=1092"),
    @("24×52=1248", "14×20=280"),
    @("20×25=500", "47×60=2820"),
    @("41×89=3649", "54×26=1404"),
    @("18×41=738", "26×20=520"),
    @("38×99=3762", "47×52=2444"),
    @("37×94=3478", "49×75=3675"),
    @("13×69=897", "99×11=1089"),
    @("68×15=1020", "66×26=1716"),
    @("45×29=1305", "96×39=3744")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
